# Anonymize "fedcore" -> "approach" in the header rows of both sheets, and
# give the merged header cells (C1/D1 on sheet 1; C1/D1/F1/G1 on sheet 2)
# the partial borders that a merged B1:D1 / E1:G1 range needs: the
# left-most cell of each merged block keeps its original full box border,
# the interior cell gets only a top+bottom border, and the right-most
# cell gets a top+bottom+right border.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("quality_comparison")
$ws2 = $wb.Worksheets.Item("computational_comparison")

# ---- Build the two border patterns once on sheet 1, then clone the
# resulting cell format onto every other cell that needs it. Re-deriving
# the patterns edge-by-edge on every single cell leaves unused entries in
# the style table, so we compute each pattern exactly once and reuse it.

# Interior-of-merge pattern: top + bottom only
$c1 = $ws1.Range("C1")
$c1.Style = "Normal"
$c1.Borders.Item(8).LineStyle = 1   # xlEdgeTop
$c1.Borders.Item(9).LineStyle = 1   # xlEdgeBottom

# Right-edge-of-merge pattern: top + bottom + right
$d1 = $ws1.Range("D1")
$d1.Style = "Normal"
$d1.Borders.Item(10).LineStyle = 1  # xlEdgeRight
$d1.Borders.Item(8).LineStyle = 1   # xlEdgeTop
$d1.Borders.Item(9).LineStyle = 1   # xlEdgeBottom

# Clone the two patterns onto sheet 2's matching cells (C1/F1 = interior,
# D1/G1 = right edge) via a formats-only paste.
$c1.Copy()
$ws2.Range("C1").PasteSpecial(-4122)  # xlPasteFormats
$ws2.Range("F1").PasteSpecial(-4122)

$d1.Copy()
$ws2.Range("D1").PasteSpecial(-4122)
$ws2.Range("G1").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---- Text anonymization: "fedcore" -> "approach" ----
$ws1.Range("C2").Value2 = "approach"
$ws2.Range("C2").Value2 = "approach"
$ws2.Range("F2").Value2 = "approach"

# ---- Drop the stray empty inline-string cell at G5 (the model_size
# "change" column is blank because E5/F5 are both 0, so there is nothing
# to show there). ----
$ws2.Range("G5").ClearContents()
